$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 30182.924133
$ws.Range("D2").Value = 786.054624

$ws.Range("B3").Value = 1782.920878
$ws.Range("D3").Value = 23.216326
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = 12748.135444
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = -1.97382
$ws.Range("H5").Value = -4.220485
$ws.Range("I5").Value = 0.272845
$ws.Range("J5").Value = 0.09804599999999999

$ws.Range("G6").Value = 3.159976
$ws.Range("H6").Value = 0.800216
$ws.Range("I6").Value = 5.519736
$ws.Range("J6").Value = 0.00501

$ws.Range("G7").Value = 5.133796
$ws.Range("H7").Value = 3.359414
$ws.Range("I7").Value = 6.908178
$ws.Range("J7").Value = 0
